# Updates the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# tracker sheet with freshly scraped values, matching the GitHub Actions
# scheduled refresh.
#
# D-column values that are unambiguous numeric strings (e.g. "587.25")
# would otherwise be auto-converted to numbers by Excel's type inference,
# dropping the trailing zero/format. They are written with a leading
# apostrophe so Excel stores them as literal text, exactly as the source
# feed provides them. D-column values that already contain multiple dots
# (e.g. "63.051.54") are never numeric-looking, so no apostrophe is
# needed there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.051.54"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "3.147.60"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'587.25"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").Value = "'137.93"
$ws.Range("E6").Value = "  -3.46%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "3.144.15"
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").Value = "'5.28"
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").Value = "'0.457"
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").Value = "'0.0000244"
$ws.Range("E13").Value = "  -2.51%  "
$ws.Range("D14").Value = "'34.09"
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("D15").Value = "3.671.03"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").Value = "3.147.37"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("D18").Value = "63.082.98"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").Value = "'6.65"
$ws.Range("E19").Value = "  -1.56%  "
$ws.Range("D20").Value = "'475.42"
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("D21").Value = "'13.99"
$ws.Range("E21").Value = "  -3.75%  "
$ws.Range("D22").Value = "'0.701"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").Value = "'7.71"
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("D24").Value = "'84.56"
$ws.Range("E24").Value = "  -2.97%  "
$ws.Range("D25").Value = "'12.95"
$ws.Range("E25").Value = "  -2.35%  "
$ws.Range("D27").Value = "'2.70"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").Value = "'7.08"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").Value = "'7.92"
$ws.Range("E29").Value = "  -3.72%  "
$ws.Range("E30").Value = "  +3.32%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "'26.74"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("E33").Value = "  -4.53%  "
$ws.Range("D34").Value = "'2.53"
$ws.Range("E34").Value = "  -3.83%  "
$ws.Range("E35").Value = "  -2.24%  "
$ws.Range("D36").Value = "'5.79"
$ws.Range("E36").Value = "  -3.00%  "
$ws.Range("D37").Value = "'52.39"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "0.0₃0697"
$ws.Range("E38").Value = "  -7.04%  "
$ws.Range("D39").Value = "'0.0388"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "'415.13"
$ws.Range("E40").Value = "  -4.61%  "
$ws.Range("D41").Value = "'2.74"
$ws.Range("E41").Value = "  -6.02%  "
$ws.Range("D42").Value = "'8.26"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").Value = "2.926.25"
$ws.Range("E43").Value = "  +2.49%  "
$ws.Range("E44").Value = "  -6.29%  "
$ws.Range("D45").Value = "'0.260"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'2.12"
$ws.Range("E47").Value = "  -3.62%  "
$ws.Range("D48").Value = "'25.37"
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "'2.23"
$ws.Range("E50").Value = "  -7.96%  "
$ws.Range("D51").Value = "'120.85"
$ws.Range("E51").Value = "  -0.35%  "

